$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (rich-text shared strings collapse to plain text; formatting unchanged visually)
$ws.Range("A8").Value = "Volume 30   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# Weekly crime statistics data refresh (rows 14-30)
# Row 14
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 12
$ws.Range("E14").Value = -16.666666666666
$ws.Range("F14").Value = 31
$ws.Range("G14").Value = 29
$ws.Range("H14").Value = 6.896551724137
$ws.Range("I14").Value = 193
$ws.Range("J14").Value = 209
$ws.Range("K14").Value = -7.655502392344
$ws.Range("L14").Value = -14.601769911504
$ws.Range("M14").Value = -14.977973568281
$ws.Range("N14").Value = -78.57935627081
# Row 15
$ws.Range("C15").Value = 21
$ws.Range("D15").Value = 38
$ws.Range("E15").Value = -44.736842105263
$ws.Range("F15").Value = 99
$ws.Range("G15").Value = 137
$ws.Range("H15").Value = -27.737226277372
$ws.Range("I15").Value = 708
$ws.Range("J15").Value = 790
$ws.Range("K15").Value = -10.379746835443
$ws.Range("L15").Value = 3.660322108345
$ws.Range("M15").Value = 16.256157635468
$ws.Range("N15").Value = -54.789272030651
# Row 16
$ws.Range("C16").Value = 346
$ws.Range("D16").Value = 376
$ws.Range("E16").Value = -7.978723404255
$ws.Range("F16").Value = 1304
$ws.Range("G16").Value = 1457
$ws.Range("H16").Value = -10.501029512697
$ws.Range("I16").Value = 7535
$ws.Range("J16").Value = 7888
$ws.Range("K16").Value = -4.475152129817
$ws.Range("L16").Value = 33.694109297374
$ws.Range("M16").Value = -12.475316529213
$ws.Range("N16").Value = -81.156847054116
# Row 17
$ws.Range("C17").Value = 637
$ws.Range("D17").Value = 583
$ws.Range("E17").Value = 9.26243567753
$ws.Range("F17").Value = 2349
$ws.Range("G17").Value = 2364
$ws.Range("H17").Value = -0.634517766497
$ws.Range("I17").Value = 12959
$ws.Range("J17").Value = 12181
$ws.Range("K17").Value = 6.386996141531
$ws.Range("L17").Value = 28.319635607485
$ws.Range("M17").Value = 62.210539491801
$ws.Range("N17").Value = -32.55789747593
# Row 18
$ws.Range("C18").Value = 257
$ws.Range("D18").Value = 315
$ws.Range("E18").Value = -18.412698412698
$ws.Range("F18").Value = 933
$ws.Range("G18").Value = 1251
$ws.Range("H18").Value = -25.419664268585
$ws.Range("I18").Value = 6687
$ws.Range("J18").Value = 7440
$ws.Range("K18").Value = -10.120967741935
$ws.Range("L18").Value = 21.383191141768
$ws.Range("M18").Value = -20.863905325443
$ws.Range("N18").Value = -85.848817031362
# Row 19
$ws.Range("C19").Value = 996
$ws.Range("D19").Value = 1025
$ws.Range("E19").Value = -2.829268292682
$ws.Range("F19").Value = 3945
$ws.Range("G19").Value = 4049
$ws.Range("H19").Value = -2.568535440849
$ws.Range("I19").Value = 23529
$ws.Range("J19").Value = 23809
$ws.Range("K19").Value = -1.176025872569
$ws.Range("L19").Value = 48.363705151648
$ws.Range("M19").Value = 37.693117977528
$ws.Range("N19").Value = -40.160223804679
# Row 20
$ws.Range("C20").Value = 353
$ws.Range("D20").Value = 280
$ws.Range("E20").Value = 26.071428571428
$ws.Range("F20").Value = 1266
$ws.Range("G20").Value = 1031
$ws.Range("H20").Value = 22.793404461687
$ws.Range("I20").Value = 7305
$ws.Range("J20").Value = 6178
$ws.Range("K20").Value = 18.242149562965
$ws.Range("L20").Value = 74.593690248566
$ws.Range("M20").Value = 51.55601659751
$ws.Range("N20").Value = -86.415115392483
# Row 21
$ws.Range("C21").Value = 2620
$ws.Range("D21").Value = 2629
$ws.Range("E21").Value = -0.342335488779
$ws.Range("F21").Value = 9927
$ws.Range("G21").Value = 10318
$ws.Range("H21").Value = -3.789494088001
$ws.Range("I21").Value = 58916
$ws.Range("J21").Value = 58495
$ws.Range("K21").Value = 0.719719634156
$ws.Range("L21").Value = 39.624608967674
$ws.Range("M21").Value = 23.275862068965
$ws.Range("N21").Value = -70.836117752466
# Row 22
$ws.Range("C22").Value = 43
$ws.Range("D22").Value = 38
$ws.Range("E22").Value = 13.157894736842
$ws.Range("F22").Value = 171
$ws.Range("G22").Value = 164
$ws.Range("H22").Value = 4.268292682926
$ws.Range("I22").Value = 1053
$ws.Range("J22").Value = 1112
$ws.Range("K22").Value = -5.305755395683
$ws.Range("L22").Value = 46.047156726768
$ws.Range("M22").Value = 5.3
# Row 23
$ws.Range("C23").Value = 115
$ws.Range("D23").Value = 136
$ws.Range("E23").Value = -15.441176470588
$ws.Range("F23").Value = 491
$ws.Range("G23").Value = 557
$ws.Range("H23").Value = -11.849192100538
$ws.Range("I23").Value = 2936
$ws.Range("J23").Value = 2841
$ws.Range("K23").Value = 3.343892995424
$ws.Range("L23").Value = 17.486994797919
$ws.Range("M23").Value = 54.445028932141
# Row 24
$ws.Range("C24").Value = 2320
$ws.Range("D24").Value = 2477
$ws.Range("E24").Value = -6.338312474767
$ws.Range("F24").Value = 9011
$ws.Range("G24").Value = 9471
$ws.Range("H24").Value = -4.856931686199
$ws.Range("I24").Value = 52237
$ws.Range("J24").Value = 53275
$ws.Range("K24").Value = -1.948381041764
$ws.Range("L24").Value = 40.33527657631
$ws.Range("M24").Value = 40.083132207026
# Row 25
$ws.Range("C25").Value = 891
$ws.Range("D25").Value = 884
$ws.Range("E25").Value = 0.791855203619
$ws.Range("F25").Value = 3680
$ws.Range("G25").Value = 3691
$ws.Range("H25").Value = -0.298022216201
$ws.Range("I25").Value = 20776
$ws.Range("J25").Value = 19832
$ws.Range("K25").Value = 4.759983864461
$ws.Range("L25").Value = 32.660749632845
$ws.Range("M25").Value = -6.029218870143
# Row 26
$ws.Range("C26").Value = 36
$ws.Range("D26").Value = 60
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 175
$ws.Range("G26").Value = 224
$ws.Range("H26").Value = -21.875
$ws.Range("I26").Value = 1182
$ws.Range("J26").Value = 1300
$ws.Range("K26").Value = -9.076923076923
$ws.Range("L26").Value = 4.509283819628
# Row 27
$ws.Range("C27").Value = 94
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = -6
$ws.Range("F27").Value = 428
$ws.Range("G27").Value = 439
$ws.Range("H27").Value = -2.50569476082
$ws.Range("I27").Value = 2514
$ws.Range("J27").Value = 2427
$ws.Range("K27").Value = 3.584672435105
$ws.Range("L27").Value = 16.442797591477
# Row 28
$ws.Range("C28").Value = 35
$ws.Range("D28").Value = 59
$ws.Range("E28").Value = -40.677966101694
$ws.Range("F28").Value = 102
$ws.Range("G28").Value = 158
$ws.Range("H28").Value = -35.443037974683
$ws.Range("I28").Value = 543
$ws.Range("J28").Value = 739
$ws.Range("K28").Value = -26.522327469553
$ws.Range("L28").Value = -33.210332103321
$ws.Range("M28").Value = -29.296875
$ws.Range("N28").Value = -80.051432770022
# Row 29
$ws.Range("C29").Value = 28
$ws.Range("D29").Value = 35
$ws.Range("E29").Value = -20
$ws.Range("F29").Value = 91
$ws.Range("G29").Value = 122
$ws.Range("H29").Value = -25.409836065573
$ws.Range("I29").Value = 464
$ws.Range("J29").Value = 616
$ws.Range("K29").Value = -24.675324675324
$ws.Range("L29").Value = -34.370579915134
$ws.Range("M29").Value = -26.465927099841
$ws.Range("N29").Value = -81.184103811841
# Row 30
$ws.Range("C30").Value = 9
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = -10
$ws.Range("F30").Value = 28
$ws.Range("G30").Value = 55
$ws.Range("H30").Value = -49.090909090909
$ws.Range("I30").Value = 229
$ws.Range("J30").Value = 337
$ws.Range("K30").Value = -32.047477744807
$ws.Range("L30").Value = -20.209059233449
